{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst marker = \"4- Seria uma equipe de seis a oito pessoas. \";\nlet target = null;\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (const p of paragraphs.items) {\n  if (p.text && p.text.indexOf(marker) === 0) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Target paragraph not found\");\n}\n\nconst results = target.search(marker, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Marker text not found in target paragraph\");\n}\n\nconst hit = results.items[0];\nhit.insertText(\"Seria feito no modelo espiral. \", \"After\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"4- Seria uma equipe de seis a oito pessoas. \"\n$find.Execute()\nif ($find.Found) {\n    $r = $find.Parent\n    $r.Collapse(0)\n    $r.InsertAfter(\"Seria feito no modelo espiral. \")\n}\n"}
